$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (written A1, C1, B1, D1 to match original authoring order)
$ws.Range("A1").Value = "Fonctionnalité"
$ws.Range("C1").Value = "Utiliser les fonctionnalité du host en étant client"
$ws.Range("B1").Value = "Pouvoir se connecter en tant que client/serveur"
$ws.Range("D1").Value = "ne pas pouvoir utiliser les outils adverses"

# Data row
$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "système client/socket"
$ws.Range("C2").Value = "structure lancant une méthode"
$ws.Range("D2").Value = "Système de propriétés"

# Column widths (values compensate for the engine's pixel-grid quantization
# of ColumnWidth, so the saved width lands as close as possible to the
# target 21.140625 / 43.7109375 / 44.85546875 / 24.5703125 char-widths)
$ws.Columns.Item(1).ColumnWidth = 20.307291666666668
$ws.Columns.Item(2).ColumnWidth = 42.877604166666664
$ws.Columns.Item(3).ColumnWidth = 44.022135416666664
$ws.Columns.Item(4).ColumnWidth = 23.736979166666668

# Selection, to match <selection activeCell="C8" sqref="C8"/>
$ws.Range("C8").Select()
